$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) "301.49"
Set-TextValue $ws.Cells.Item(2, 5) "-1.26%"
Set-TextValue $ws.Cells.Item(2, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(2, 7) "1"

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) "37.26"
Set-TextValue $ws.Cells.Item(3, 5) "5.40%"
Set-TextValue $ws.Cells.Item(3, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(3, 7) "1"

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) "4.998"
Set-TextValue $ws.Cells.Item(4, 5) "-4.49%"
Set-TextValue $ws.Cells.Item(4, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(4, 7) "1"

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) "0.07716"
Set-TextValue $ws.Cells.Item(5, 5) "-1.27%"
Set-TextValue $ws.Cells.Item(5, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(5, 7) "1"

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) "2.229"
Set-TextValue $ws.Cells.Item(6, 5) "-6.95%"
Set-TextValue $ws.Cells.Item(6, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(6, 7) "1"

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) "8.001"
Set-TextValue $ws.Cells.Item(7, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(7, 7) "1"

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) "4.002"
Set-TextValue $ws.Cells.Item(8, 5) "1.94%"
Set-TextValue $ws.Cells.Item(8, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(8, 7) "1"

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) "0.9228"
Set-TextValue $ws.Cells.Item(9, 5) "-1.36%"
Set-TextValue $ws.Cells.Item(9, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(9, 7) "1"

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) "0.09261"
Set-TextValue $ws.Cells.Item(10, 5) "-5.61%"
Set-TextValue $ws.Cells.Item(10, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(10, 7) "1"

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) "0.1820"
Set-TextValue $ws.Cells.Item(11, 5) "1.61%"
Set-TextValue $ws.Cells.Item(11, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(11, 7) "1"

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) "0.08428"
Set-TextValue $ws.Cells.Item(12, 5) "-1.78%"
Set-TextValue $ws.Cells.Item(12, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(12, 7) "1"

# Row 13
Set-TextValue $ws.Cells.Item(13, 4) "0.03603"
Set-TextValue $ws.Cells.Item(13, 5) "8.17%"
Set-TextValue $ws.Cells.Item(13, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(13, 7) "1"

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) "0.09946"
Set-TextValue $ws.Cells.Item(14, 5) "0.50%"
Set-TextValue $ws.Cells.Item(14, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(14, 7) "1"

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) "0.001478"
Set-TextValue $ws.Cells.Item(15, 5) "-0.75%"
Set-TextValue $ws.Cells.Item(15, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(15, 7) "1"

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) "0.005717"
Set-TextValue $ws.Cells.Item(16, 5) "-1.73%"
Set-TextValue $ws.Cells.Item(16, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(16, 7) "1"

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) "3.470"
Set-TextValue $ws.Cells.Item(17, 5) "0.20%"
Set-TextValue $ws.Cells.Item(17, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(17, 7) "1"

# Row 18
Set-TextValue $ws.Cells.Item(18, 4) "2.133"
Set-TextValue $ws.Cells.Item(18, 5) "-4.06%"
Set-TextValue $ws.Cells.Item(18, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(18, 7) "1"

# Row 19
Set-TextValue $ws.Cells.Item(19, 5) "2.87%"
Set-TextValue $ws.Cells.Item(19, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(19, 7) "1"

# Row 20
Set-TextValue $ws.Cells.Item(20, 4) "0.1303"
Set-TextValue $ws.Cells.Item(20, 5) "-1.47%"
Set-TextValue $ws.Cells.Item(20, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(20, 7) "1"

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) "4.555"
Set-TextValue $ws.Cells.Item(21, 5) "4.07%"
Set-TextValue $ws.Cells.Item(21, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(21, 7) "1"

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) "0.2238"
Set-TextValue $ws.Cells.Item(22, 5) "-2.83%"
Set-TextValue $ws.Cells.Item(22, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(22, 7) "1"

# Row 23
Set-TextValue $ws.Cells.Item(23, 5) "0.74%"
Set-TextValue $ws.Cells.Item(23, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(23, 7) "1"

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) "0.001231"
Set-TextValue $ws.Cells.Item(24, 5) "2.19%"
Set-TextValue $ws.Cells.Item(24, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(24, 7) "1"

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) "0.004453"
Set-TextValue $ws.Cells.Item(25, 5) "0.85%"
Set-TextValue $ws.Cells.Item(25, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(25, 7) "1"

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) "0.0001302"
Set-TextValue $ws.Cells.Item(26, 5) "0.02%"
Set-TextValue $ws.Cells.Item(26, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(26, 7) "1"

# Row 27
Set-TextValue $ws.Cells.Item(27, 5) "-20.55%"
Set-TextValue $ws.Cells.Item(27, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(27, 7) "1"

# Row 28
Set-TextValue $ws.Cells.Item(28, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(28, 7) "1"

# Row 29
Set-TextValue $ws.Cells.Item(29, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(29, 7) "1"

# Row 30
Set-TextValue $ws.Cells.Item(30, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(30, 7) "1"

# Row 31
Set-TextValue $ws.Cells.Item(31, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(31, 7) "1"

# Row 32
Set-TextValue $ws.Cells.Item(32, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(32, 7) "1"

# Row 33
Set-TextValue $ws.Cells.Item(33, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(33, 7) "1"

# Row 34
Set-TextValue $ws.Cells.Item(34, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(34, 7) "1"

# Row 35
Set-TextValue $ws.Cells.Item(35, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(35, 7) "1"

# Row 36
Set-TextValue $ws.Cells.Item(36, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(36, 7) "1"

# Row 37
Set-TextValue $ws.Cells.Item(37, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(37, 7) "1"

# Row 38
Set-TextValue $ws.Cells.Item(38, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(38, 7) "1"

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) "0.01743"
Set-TextValue $ws.Cells.Item(39, 5) "-1.27%"
Set-TextValue $ws.Cells.Item(39, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(39, 7) "1"

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) "0.04700"
Set-TextValue $ws.Cells.Item(40, 5) "-1.87%"
Set-TextValue $ws.Cells.Item(40, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(40, 7) "1"

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) "0.007860"
Set-TextValue $ws.Cells.Item(41, 5) "1.08%"
Set-TextValue $ws.Cells.Item(41, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(41, 7) "1"

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) "0.1392"
Set-TextValue $ws.Cells.Item(42, 5) "-1.75%"
Set-TextValue $ws.Cells.Item(42, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(42, 7) "1"

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) "0.007686"
Set-TextValue $ws.Cells.Item(43, 5) "-21.61%"
Set-TextValue $ws.Cells.Item(43, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(43, 7) "1"

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) "0.002226"
Set-TextValue $ws.Cells.Item(44, 5) "6.28%"
Set-TextValue $ws.Cells.Item(44, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(44, 7) "1"

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) "0.009027"
Set-TextValue $ws.Cells.Item(45, 5) "-10.72%"
Set-TextValue $ws.Cells.Item(45, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(45, 7) "1"

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) "0.00006204"
Set-TextValue $ws.Cells.Item(46, 5) "1.14%"
Set-TextValue $ws.Cells.Item(46, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(46, 7) "1"

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) "0.00000000751"
Set-TextValue $ws.Cells.Item(47, 5) "0.07%"
Set-TextValue $ws.Cells.Item(47, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(47, 7) "1"

# Row 48
Set-TextValue $ws.Cells.Item(48, 5) "19.83%"
Set-TextValue $ws.Cells.Item(48, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(48, 7) "1"

# Row 49
Set-TextValue $ws.Cells.Item(49, 4) "0.002699"
Set-TextValue $ws.Cells.Item(49, 5) "34.77%"
Set-TextValue $ws.Cells.Item(49, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(49, 7) "1"

# Row 50
Set-TextValue $ws.Cells.Item(50, 4) "0.00002104"
Set-TextValue $ws.Cells.Item(50, 5) "0.07%"
Set-TextValue $ws.Cells.Item(50, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(50, 7) "1"

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) "0.0002004"
Set-TextValue $ws.Cells.Item(51, 5) "0.07%"
Set-TextValue $ws.Cells.Item(51, 6) "22-1-2023"
Set-TextValue $ws.Cells.Item(51, 7) "1"
